$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.56
$ws.Range("H2").Value = 4.1
$ws.Range("I2").Value = 4.3
$ws.Range("K2").Value = 2.76
$ws.Range("L2").Value = 1.65
$ws.Range("P2").Value = 1.44
$ws.Range("R2").Value = 1.17
$ws.Range("S2").Value = 6.4
$ws.Range("T2").Value = 2.34
$ws.Range("U2").Value = 1.57
$ws.Range("V2").Value = 1.3
$ws.Range("X2").Value = 6.6
$ws.Range("Y2").Value = 16
$ws.Range("Z2").Value = 1000
$ws.Range("AA2").Value = 910
$ws.Range("AB2").Value = 7
$ws.Range("AD2").Value = 29
$ws.Range("AF2").Value = 22
$ws.Range("AG2").Value = 11
$ws.Range("AH2").Value = 34
$ws.Range("AI2").Value = 170
$ws.Range("AJ2").Value = 65
$ws.Range("AK2").Value = 38
$ws.Range("AL2").Value = 120
$ws.Range("AN2").Value = 50
$ws.Range("AO2").Value = 600
# Row 3
$ws.Range("H3").Value = 2.24
$ws.Range("J3").Value = 3.65
$ws.Range("N3").Value = 4.6
$ws.Range("T3").Value = 1.64
$ws.Range("V3").Value = 1.75
$ws.Range("X3").Value = 18
$ws.Range("Z3").Value = 16
$ws.Range("AA3").Value = 30
$ws.Range("AC3").Value = 8.4
$ws.Range("AD3").Value = 11.5
$ws.Range("AE3").Value = 23
$ws.Range("AI3").Value = 32
$ws.Range("AJ3").Value = 60
$ws.Range("AK3").Value = 36
$ws.Range("AL3").Value = 42
$ws.Range("AM3").Value = 75
$ws.Range("AN3").Value = 32
# Row 4
$ws.Range("F4").Value = 2.06
$ws.Range("G4").Value = 2.22
$ws.Range("H4").Value = 3.7
$ws.Range("I4").Value = 4.3
$ws.Range("J4").Value = 3.5
$ws.Range("K4").Value = 3.75
$ws.Range("L4").Value = 1.45
$ws.Range("N4").Value = 3.5
$ws.Range("O4").Value = 1.35
$ws.Range("P4").Value = 1.83
$ws.Range("Q4").Value = 2.06
$ws.Range("S4").Value = 3.7
$ws.Range("T4").Value = 1.84
$ws.Range("U4").Value = 1.96
$ws.Range("V4").Value = 1.32
$ws.Range("W4").Value = 1.82
$ws.Range("X4").Value = 18
$ws.Range("Y4").Value = 15
$ws.Range("Z4").Value = 1000
$ws.Range("AA4").Value = 1000
$ws.Range("AB4").Value = 9.6
$ws.Range("AC4").Value = 8.4
$ws.Range("AE4").Value = 60
$ws.Range("AF4").Value = 14.5
$ws.Range("AG4").Value = 11
$ws.Range("AI4").Value = 65
$ws.Range("AJ4").Value = 200
$ws.Range("AK4").Value = 25
$ws.Range("AL4").Value = 55
$ws.Range("AM4").Value = 140
$ws.Range("AN4").Value = 19.5
$ws.Range("AO4").Value = 1000
# Row 5
$ws.Range("F5").Value = 1.88
$ws.Range("G5").Value = 1.96
$ws.Range("H5").Value = 5.7
$ws.Range("I5").Value = 6.8
$ws.Range("K5").Value = 3.4
$ws.Range("P5").Value = 1.45
$ws.Range("Q5").Value = 2.96
$ws.Range("R5").Value = 1.15
$ws.Range("S5").Value = 6.4
$ws.Range("T5").Value = 2.42
$ws.Range("U5").Value = 1.57
$ws.Range("V5").Value = 1.19
$ws.Range("W5").Value = 2.04
$ws.Range("Z5").Value = 50
$ws.Range("AA5").Value = 210
$ws.Range("AB5").Value = 5.8
$ws.Range("AC5").Value = 7.6
$ws.Range("AD5").Value = 28
$ws.Range("AF5").Value = 9.800000000000001
$ws.Range("AH5").Value = 40
$ws.Range("AN5").Value = 30
$ws.Range("AO5").Value = 280
# Row 6
$ws.Range("F6").Value = 4.1
$ws.Range("G6").Value = 4.2
$ws.Range("I6").Value = 2.14
$ws.Range("P6").Value = 1.64
$ws.Range("S6").Value = 5
$ws.Range("U6").Value = 1.83
$ws.Range("V6").Value = 1.87
$ws.Range("W6").Value = 1.31
$ws.Range("AC6").Value = 7.4
$ws.Range("AN6").Value = 160
$ws.Range("AO6").Value = 25
# Row 7
$ws.Range("G7").Value = 1.43
# Row 8
$ws.Range("F8").Value = 2.56
$ws.Range("G8").Value = 2.74
$ws.Range("H8").Value = 3.05
$ws.Range("I8").Value = 3.3
$ws.Range("J8").Value = 3.1
$ws.Range("K8").Value = 3.3
$ws.Range("P8").Value = 1.84
$ws.Range("R8").Value = 1.33
$ws.Range("S8").Value = 3.65
$ws.Range("T8").Value = 1.74
$ws.Range("U8").Value = 2.08
$ws.Range("V8").Value = 1.43
$ws.Range("W8").Value = 1.57
$ws.Range("Y8").Value = 12.5
$ws.Range("Z8").Value = 21
$ws.Range("AB8").Value = 12
$ws.Range("AE8").Value = 38
$ws.Range("AF8").Value = 30
$ws.Range("AJ8").Value = 44
$ws.Range("AK8").Value = 29
$ws.Range("AN8").Value = 55
$ws.Range("AO8").Value = 46
# Row 9
$ws.Range("Q9").Value = 1.73
$ws.Range("R9").Value = 1.51
$ws.Range("W9").Value = 2.8
$ws.Range("AN9").Value = 7
$ws.Range("AO9").Value = 110
# Row 11
$ws.Range("F11").Value = 2.48
# Row 12
$ws.Range("F12").Value = 2
$ws.Range("K12").Value = 3.5
$ws.Range("S12").Value = 4
# Row 13
$ws.Range("F13").Value = 2.26
$ws.Range("G13").Value = 2.28
$ws.Range("H13").Value = 3.85
$ws.Range("I13").Value = 4
$ws.Range("J13").Value = 3.25
$ws.Range("K13").Value = 3.3
$ws.Range("N13").Value = 2.78
$ws.Range("V13").Value = 1.33
$ws.Range("W13").Value = 1.78
$ws.Range("Z13").Value = 25
$ws.Range("AA13").Value = 95
$ws.Range("AD13").Value = 17
$ws.Range("AE13").Value = 65
$ws.Range("AI13").Value = 90
$ws.Range("AO13").Value = 90
# Row 14
$ws.Range("F14").Value = 1.61
$ws.Range("H14").Value = 5.9
$ws.Range("I14").Value = 6.8
$ws.Range("J14").Value = 3.95
$ws.Range("Q14").Value = 1.91
$ws.Range("S14").Value = 3.2
$ws.Range("U14").Value = 1.97
$ws.Range("Y14").Value = 1000
# Row 15
$ws.Range("V15").Value = 1.25
# Row 16
$ws.Range("F16").Value = 2.08
$ws.Range("Z16").Value = 980
$ws.Range("AA16").Value = 75
$ws.Range("AI16").Value = 48
$ws.Range("AK16").Value = 21
$ws.Range("AN16").Value = 13.5
